$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 146.17053
$ws.Range("H2").Value = 438.51159
$ws.Range("I2").Value = 0.4047435297111188
$ws.Range("J2").Value = 0.4047435297111188
$ws.Range("M2").Value = 27.681071
$ws.Range("N2").Value = 83.04321300000001
$ws.Range("O2").Value = 0.05045805550111082
$ws.Range("P2").Value = 0.05045805550111081
$ws.Range("Q2").Value = 4046.15681903763
$ws.Range("R2").Value = 36415.41137133867
$ws.Range("S2").Value = 0.02042257148587913
$ws.Range("T2").Value = 0.02042257148587913
$ws.Range("G3").Value = 146.17053
$ws.Range("H3").Value = 438.51159
$ws.Range("I3").Value = 0.4047435297111188
$ws.Range("J3").Value = 0.4047435297111188
$ws.Range("O3").Value = 0.0001771869602491167
$ws.Range("P3").Value = 0.0001771869602491166
$ws.Range("Q3").Value = 14.20836019812
$ws.Range("R3").Value = 127.87524178308
$ws.Range("S3").Value = 0.00007171527571001116
$ws.Range("T3").Value = 0.00007171527571001115
$ws.Range("G4").Value = 146.17053
$ws.Range("H4").Value = 438.51159
$ws.Range("I4").Value = 0.4047435297111188
$ws.Range("J4").Value = 0.4047435297111188
$ws.Range("M4").Value = 272.2666776666667
$ws.Range("N4").Value = 816.800033
$ws.Range("O4").Value = 0.4962975288350554
$ws.Range("P4").Value = 0.4962975288350553
$ws.Range("Q4").Value = 39797.36457587583
$ws.Range("R4").Value = 358176.2811828824
$ws.Range("S4").Value = 0.2008732136076061
$ws.Range("T4").Value = 0.200873213607606
$ws.Range("G5").Value = 146.17053
$ws.Range("H5").Value = 438.51159
$ws.Range("I5").Value = 0.4047435297111188
$ws.Range("J5").Value = 0.4047435297111188
$ws.Range("M5").Value = 11.73516533333333
$ws.Range("N5").Value = 35.205496
$ws.Range("O5").Value = 0.02139128300722342
$ws.Range("P5").Value = 0.02139128300722341
$ws.Range("Q5").Value = 1715.33533641096
$ws.Range("R5").Value = 15438.01802769864
$ws.Range("S5").Value = 0.008657983389393081
$ws.Range("T5").Value = 0.00865798338939308
$ws.Range("G6").Value = 146.17053
$ws.Range("H6").Value = 438.51159
$ws.Range("I6").Value = 0.4047435297111188
$ws.Range("J6").Value = 0.4047435297111188
$ws.Range("M6").Value = 236.8155566666667
$ws.Range("N6").Value = 710.44667
$ws.Range("O6").Value = 0.4316759456963613
$ws.Range("P6").Value = 0.4316759456963613
$ws.Range("Q6").Value = 34615.4554302117
$ws.Range("R6").Value = 311539.0988719053
$ws.Range("S6").Value = 0.1747180459525305
$ws.Range("T6").Value = 0.1747180459525305
$ws.Range("H7").Value = 632.3552549999999
$ws.Range("I7").Value = 0.5836600531814327
$ws.Range("J7").Value = 0.5836600531814327
$ws.Range("M7").Value = 27.681071
$ws.Range("N7").Value = 83.04321300000001
$ws.Range("O7").Value = 0.05045805550111082
$ws.Range("P7").Value = 0.05045805550111081
$ws.Range("Q7").Value = 5834.756903626035
$ws.Range("R7").Value = 52512.81213263432
$ws.Range("S7").Value = 0.02945035135721002
$ws.Range("T7").Value = 0.02945035135721002
$ws.Range("H8").Value = 632.3552549999999
$ws.Range("I8").Value = 0.5836600531814327
$ws.Range("J8").Value = 0.5836600531814327
$ws.Range("O8").Value = 0.0001771869602491167
$ws.Range("P8").Value = 0.0001771869602491166
$ws.Range("S8").Value = 0.0001034169506420558
$ws.Range("T8").Value = 0.0001034169506420558
$ws.Range("H9").Value = 632.3552549999999
$ws.Range("I9").Value = 0.5836600531814327
$ws.Range("J9").Value = 0.5836600531814327
$ws.Range("M9").Value = 272.2666776666667
$ws.Range("N9").Value = 816.800033
$ws.Range("O9").Value = 0.4962975288350554
$ws.Range("P9").Value = 0.4962975288350553
$ws.Range("Q9").Value = 57389.75479463593
$ws.Range("R9").Value = 516507.7931517233
$ws.Range("S9").Value = 0.289669042073682
$ws.Range("T9").Value = 0.289669042073682
$ws.Range("H10").Value = 632.3552549999999
$ws.Range("I10").Value = 0.5836600531814327
$ws.Range("J10").Value = 0.5836600531814327
$ws.Range("M10").Value = 11.73516533333333
$ws.Range("N10").Value = 35.205496
$ws.Range("O10").Value = 0.02139128300722342
$ws.Range("P10").Value = 0.02139128300722341
$ws.Range("Q10").Value = 2473.59782227572
$ws.Range("R10").Value = 22262.38040048148
$ws.Range("S10").Value = 0.0124852373776151
$ws.Range("T10").Value = 0.01248523737761509
$ws.Range("H11").Value = 632.3552549999999
$ws.Range("I11").Value = 0.5836600531814327
$ws.Range("J11").Value = 0.5836600531814327
$ws.Range("M11").Value = 236.8155566666667
$ws.Range("N11").Value = 710.44667
$ws.Range("O11").Value = 0.4316759456963613
$ws.Range("P11").Value = 0.4316759456963613
$ws.Range("Q11").Value = 49917.18724130565
$ws.Range("R11").Value = 449254.6851717508
$ws.Range("S11").Value = 0.2519520054222835
$ws.Range("T11").Value = 0.2519520054222835
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 0.5494936666666667
$ws.Range("H12").Value = 1.648481
$ws.Range("I12").Value = 0.00152153793381314
$ws.Range("J12").Value = 0.00152153793381314
$ws.Range("M12").Value = 27.681071
$ws.Range("N12").Value = 83.04321300000001
$ws.Range("O12").Value = 0.05045805550111082
$ws.Range("P12").Value = 0.05045805550111081
$ws.Range("Q12").Value = 15.21057320105033
$ws.Range("R12").Value = 136.895158809453
$ws.Range("S12").Value = 0.00007677384551138891
$ws.Range("T12").Value = 0.0000767738455113889
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 0.5494936666666667
$ws.Range("H13").Value = 1.648481
$ws.Range("I13").Value = 0.00152153793381314
$ws.Range("J13").Value = 0.00152153793381314
$ws.Range("O13").Value = 0.0001771869602491167
$ws.Range("P13").Value = 0.0001771869602491166
$ws.Range("Q13").Value = 0.05341298237466666
$ws.Range("R13").Value = 0.480716841372
$ws.Range("S13").Value = 0.000000269596681396072
$ws.Range("T13").Value = 0.0000002695966813960719
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.5494936666666667
$ws.Range("H14").Value = 1.648481
$ws.Range("I14").Value = 0.00152153793381314
$ws.Range("J14").Value = 0.00152153793381314
$ws.Range("M14").Value = 272.2666776666667
$ws.Range("N14").Value = 816.800033
$ws.Range("O14").Value = 0.4962975288350554
$ws.Range("P14").Value = 0.4962975288350553
$ws.Range("Q14").Value = 149.6088150222081
$ws.Range("R14").Value = 1346.479335199873
$ws.Range("S14").Value = 0.0007551355165802576
$ws.Range("T14").Value = 0.0007551355165802575
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.5494936666666667
$ws.Range("H15").Value = 1.648481
$ws.Range("I15").Value = 0.00152153793381314
$ws.Range("J15").Value = 0.00152153793381314
$ws.Range("M15").Value = 11.73516533333333
$ws.Range("N15").Value = 35.205496
$ws.Range("O15").Value = 0.02139128300722342
$ws.Range("P15").Value = 0.02139128300722341
$ws.Range("Q15").Value = 6.448399027952888
$ws.Range("R15").Value = 58.03559125157599
$ws.Range("S15").Value = 0.00003254764854842286
$ws.Range("T15").Value = 0.00003254764854842285
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.5494936666666667
$ws.Range("H16").Value = 1.648481
$ws.Range("I16").Value = 0.00152153793381314
$ws.Range("J16").Value = 0.00152153793381314
$ws.Range("M16").Value = 236.8155566666667
$ws.Range("N16").Value = 710.44667
$ws.Range("O16").Value = 0.4316759456963613
$ws.Range("P16").Value = 0.4316759456963613
$ws.Range("Q16").Value = 130.1286485564744
$ws.Range("R16").Value = 1171.15783700827
$ws.Range("S16").Value = 0.0006568113264916748
$ws.Range("T16").Value = 0.0006568113264916748
$ws.Range("G17").Value = 3.410044
$ws.Range("H17").Value = 10.230132
$ws.Range("I17").Value = 0.009442349596941478
$ws.Range("J17").Value = 0.009442349596941478
$ws.Range("M17").Value = 27.681071
$ws.Range("N17").Value = 83.04321300000001
$ws.Range("O17").Value = 0.05045805550111082
$ws.Range("P17").Value = 0.05045805550111081
$ws.Range("Q17").Value = 94.393670077124
$ws.Range("R17").Value = 849.543030694116
$ws.Range("S17").Value = 0.0004764426000233645
$ws.Range("T17").Value = 0.0004764426000233644
$ws.Range("G18").Value = 3.410044
$ws.Range("H18").Value = 10.230132
$ws.Range("I18").Value = 0.009442349596941478
$ws.Range("J18").Value = 0.009442349596941478
$ws.Range("O18").Value = 0.0001771869602491167
$ws.Range("P18").Value = 0.0001771869602491166
$ws.Range("Q18").Value = 0.331469916976
$ws.Range("R18").Value = 2.983229252784
$ws.Range("S18").Value = 0.000001673061222691532
$ws.Range("T18").Value = 0.000001673061222691532
$ws.Range("G19").Value = 3.410044
$ws.Range("H19").Value = 10.230132
$ws.Range("I19").Value = 0.009442349596941478
$ws.Range("J19").Value = 0.009442349596941478
$ws.Range("M19").Value = 272.2666776666667
$ws.Range("N19").Value = 816.800033
$ws.Range("O19").Value = 0.4962975288350554
$ws.Range("P19").Value = 0.4962975288350553
$ws.Range("Q19").Value = 928.4413505771506
$ws.Range("R19").Value = 8355.972155194355
$ws.Range("S19").Value = 0.004686214771358737
$ws.Range("T19").Value = 0.004686214771358737
$ws.Range("G20").Value = 3.410044
$ws.Range("H20").Value = 10.230132
$ws.Range("I20").Value = 0.009442349596941478
$ws.Range("J20").Value = 0.009442349596941478
$ws.Range("M20").Value = 11.73516533333333
$ws.Range("N20").Value = 35.205496
$ws.Range("O20").Value = 0.02139128300722342
$ws.Range("P20").Value = 0.02139128300722341
$ws.Range("Q20").Value = 40.01743013394132
$ws.Range("R20").Value = 360.1568712054719
$ws.Range("S20").Value = 0.0002019839724813171
$ws.Range("T20").Value = 0.0002019839724813171
$ws.Range("G21").Value = 3.410044
$ws.Range("H21").Value = 10.230132
$ws.Range("I21").Value = 0.009442349596941478
$ws.Range("J21").Value = 0.009442349596941478
$ws.Range("M21").Value = 236.8155566666667
$ws.Range("N21").Value = 710.44667
$ws.Range("O21").Value = 0.4316759456963613
$ws.Range("P21").Value = 0.4316759456963613
$ws.Range("Q21").Value = 807.5514681178266
$ws.Range("R21").Value = 7267.96321306044
$ws.Range("S21").Value = 0.004076035191855368
$ws.Range("T21").Value = 0.004076035191855368
$ws.Range("E22").Value = 3.0
$ws.Range("F22").Value = 1.0
$ws.Range("G22").Value = 0.228434
$ws.Range("H22").Value = 0.6853020000000001
$ws.Range("I22").Value = 0.0006325295766939459
$ws.Range("J22").Value = 0.0006325295766939459
$ws.Range("M22").Value = 27.681071
$ws.Range("N22").Value = 83.04321300000001
$ws.Range("O22").Value = 0.05045805550111082
$ws.Range("P22").Value = 0.05045805550111081
$ws.Range("Q22").Value = 6.323297772814001
$ws.Range("R22").Value = 56.90967995532601
$ws.Range("S22").Value = 0.00003191621248691726
$ws.Range("T22").Value = 0.00003191621248691725
$ws.Range("E23").Value = 3.0
$ws.Range("F23").Value = 1.0
$ws.Range("G23").Value = 0.228434
$ws.Range("H23").Value = 0.6853020000000001
$ws.Range("I23").Value = 0.0006325295766939459
$ws.Range("J23").Value = 0.0006325295766939459
$ws.Range("O23").Value = 0.0001771869602491167
$ws.Range("P23").Value = 0.0001771869602491166
$ws.Range("Q23").Value = 0.022204698536
$ws.Range("R23").Value = 0.199842286824
$ws.Range("S23").Value = 0.0000001120759929620608
$ws.Range("T23").Value = 0.0000001120759929620608
$ws.Range("E24").Value = 3.0
$ws.Range("F24").Value = 1.0
$ws.Range("G24").Value = 0.228434
$ws.Range("H24").Value = 0.6853020000000001
$ws.Range("I24").Value = 0.0006325295766939459
$ws.Range("J24").Value = 0.0006325295766939459
$ws.Range("M24").Value = 272.2666776666667
$ws.Range("N24").Value = 816.800033
$ws.Range("O24").Value = 0.4962975288350554
$ws.Range("P24").Value = 0.4962975288350553
$ws.Range("Q24").Value = 62.19496624610734
$ws.Range("R24").Value = 559.7546962149661
$ws.Range("S24").Value = 0.000313922865828289
$ws.Range("T24").Value = 0.000313922865828289
$ws.Range("E25").Value = 3.0
$ws.Range("F25").Value = 1.0
$ws.Range("G25").Value = 0.228434
$ws.Range("H25").Value = 0.6853020000000001
$ws.Range("I25").Value = 0.0006325295766939459
$ws.Range("J25").Value = 0.0006325295766939459
$ws.Range("M25").Value = 11.73516533333333
$ws.Range("N25").Value = 35.205496
$ws.Range("O25").Value = 0.02139128300722342
$ws.Range("P25").Value = 0.02139128300722341
$ws.Range("Q25").Value = 2.680710757754667
$ws.Range("R25").Value = 24.126396819792
$ws.Range("S25").Value = 0.00001353061918549943
$ws.Range("T25").Value = 0.00001353061918549942
$ws.Range("E26").Value = 3.0
$ws.Range("F26").Value = 1.0
$ws.Range("G26").Value = 0.228434
$ws.Range("H26").Value = 0.6853020000000001
$ws.Range("I26").Value = 0.0006325295766939459
$ws.Range("J26").Value = 0.0006325295766939459
$ws.Range("M26").Value = 236.8155566666667
$ws.Range("N26").Value = 710.44667
$ws.Range("O26").Value = 0.4316759456963613
$ws.Range("P26").Value = 0.4316759456963613
$ws.Range("Q26").Value = 54.09672487159335
$ws.Range("R26").Value = 486.8705238443401
$ws.Range("S26").Value = 0.0002730478032002782
$ws.Range("T26").Value = 0.0002730478032002782
